$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for team record columns, right after the existing data (AC)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting (bold, centered, bordered) used by the rest of row 1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Every player row gets the same team record: 95 wins, 66 losses, 1 tie
$ws.Range("AD2:AD40").Value = 95
$ws.Range("AE2:AE40").Value = 66
$ws.Range("AF2:AF40").Value = 1
